$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 9: convert from 'open' style (4/5) to 'closing' style (6/7), add empty A9 ---
$ws.Range("A3:E3").Copy()
$ws.Range("A9:E9").PasteSpecial(-4122)

# --- Formatting for new rows 10-32, copied from existing template rows ---
# Type A (open / first-of-group): copy A2:E2 (style 4/5)
# Type B (closing / last-of-group): copy A3:E3 (style 6/7), includes empty A
# Type N (middle, no A cell at all): copy B2:E2 only (style 4/5)

$ws.Range("A2:E2").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)
$ws.Range("A3:E3").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)
$ws.Range("A3:E3").Copy()
$ws.Range("A13:E13").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)
$ws.Range("B2:E2").Copy()
$ws.Range("B15:E15").PasteSpecial(-4122)
$ws.Range("B2:E2").Copy()
$ws.Range("B16:E16").PasteSpecial(-4122)
$ws.Range("A3:E3").Copy()
$ws.Range("A17:E17").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A18:E18").PasteSpecial(-4122)
$ws.Range("A3:E3").Copy()
$ws.Range("A19:E19").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A20:E20").PasteSpecial(-4122)
$ws.Range("A3:E3").Copy()
$ws.Range("A21:E21").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A22:E22").PasteSpecial(-4122)
$ws.Range("B2:E2").Copy()
$ws.Range("B23:E23").PasteSpecial(-4122)
$ws.Range("A3:E3").Copy()
$ws.Range("A24:E24").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A25:E25").PasteSpecial(-4122)
$ws.Range("A3:E3").Copy()
$ws.Range("A26:E26").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A27:E27").PasteSpecial(-4122)
$ws.Range("A3:E3").Copy()
$ws.Range("A28:E28").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A29:E29").PasteSpecial(-4122)
$ws.Range("A3:E3").Copy()
$ws.Range("A30:E30").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A31:E31").PasteSpecial(-4122)
$ws.Range("B2:E2").Copy()
$ws.Range("B32:E32").PasteSpecial(-4122)

# --- Row heights ---
$ws.Rows(9).RowHeight = 21.6
$ws.Rows(10).RowHeight = 43.2
$ws.Rows(11).RowHeight = 21.6
$ws.Rows(12).RowHeight = 43.2
$ws.Rows(14).RowHeight = 43.2
$ws.Rows(15).RowHeight = 21.6
$ws.Rows(17).RowHeight = 21.6
$ws.Rows(18).RowHeight = 43.2
$ws.Rows(19).RowHeight = 21.6
$ws.Rows(20).RowHeight = 43.2
$ws.Rows(21).RowHeight = 21.6
$ws.Rows(22).RowHeight = 43.2
$ws.Rows(23).RowHeight = 21.6
$ws.Rows(24).RowHeight = 31.8
$ws.Rows(25).RowHeight = 43.2
$ws.Rows(26).RowHeight = 21.6
$ws.Rows(27).RowHeight = 43.2
$ws.Rows(29).RowHeight = 43.2
$ws.Rows(30).RowHeight = 21.6
$ws.Rows(31).RowHeight = 43.2
$ws.Rows(32).RowHeight = 21.6

# --- Cell values/content ---
$ws.Range("B9").Value = 378
$ws.Range("C9").Value = " We\'ll leave when you\'re ready."
$ws.Range("D9").Value = " Мы уйдём сразу же, как только\nты подготовишься."
$ws.Range("E9").Value = " Íú ôêäæí òñàèô çå, ëàë óïìûëï\nóú ðïäãïóïâéšûòÿ."
$ws.Range("A10").Value = "SCRIPT/D24P11A/um2001.ssb"
$ws.Range("B10").Value = 353
$ws.Range("C10").Value = " There should be a Time Gear at\nthe deepest part of this forest."
$ws.Range("D10").Value = " В дебрях этого леса должна\nнаходиться Шестерня Времени."
$ws.Range("E10").Value = " Â äåáñÿö üóïãï ìåòà äïìçîà\nîàöïäéóûòÿ Šåòóåñîÿ Âñåíåîé."
$ws.Range("B11").Value = 356
$ws.Range("C11").Value = " We\'re leaving as soon as\nyou\'re ready."
$ws.Range("D11").Value = " Мы пойдём сразу же, как ты\nподготовишься."
$ws.Range("E11").Value = " Íú ðïêäæí òñàèô çå, ëàë óú\nðïäãïóïâéšûòÿ."
$ws.Range("A12").Value = "SCRIPT/D24P11A/um2002.ssb"
$ws.Range("B12").Value = 331
$ws.Range("C12").Value = " It feels different from before...\nSomething has changed here…"
$ws.Range("D12").Value = " Здесь всё совсем иное...\nЧто-то изменилось..."
$ws.Range("E12").Value = " Èäåòû âòæ òïâòåí éîïå...\nŒóï-óï éèíåîéìïòû..."
$ws.Range("B13").Value = 334
$ws.Range("C13").Value = " Hurry. We have to go."
$ws.Range("D13").Value = " Скорее. Нам нужно идти."
$ws.Range("E13").Value = " Òëïñåå. Îàí îôçîï éäóé."
$ws.Range("A14").Value = "SCRIPT/D27P11A/um2402.ssb"
$ws.Range("B14").Value = 280
$ws.Range("C14").Value = " We have to go to [CS:P]Temporal\nTower[CR] and set the Time Gears there."
$ws.Range("D14").Value = " Нам нужно попасть в [CS:P]Темпоральную\nБашню[CR] и разместить там Шестерни Времени."
$ws.Range("E14").Value = " Îàí îôçîï ðïðàòóû â [CS:P]Óåíðïñàìûîôý\nÁàšîý[CR] é ñàèíåòóéóû óàí Šåòóåñîé Âñåíåîé."
$ws.Range("B15").Value = 283
$ws.Range("C15").Value = " That will put a stop to the\ndestruction of time."
$ws.Range("D15").Value = " Это положит конец разрушению\nвремени."
$ws.Range("E15").Value = " Üóï ðïìïçéó ëïîåø ñàèñôšåîéý\nâñåíåîé."
$ws.Range("B16").Value = 286
$ws.Range("C16").Value = " We don\'t have far to go,\n[hero]!"
$ws.Range("D16").Value = " Нам не так далеко нужно пройти,\n[hero]!"
$ws.Range("E16").Value = " Îàí îå óàë äàìåëï îôçîï ðñïêóé,\n[hero]!"
$ws.Range("B17").Value = 289
$ws.Range("C17").Value = " First, we need to reach the\nRainbow Stoneship!"
$ws.Range("D17").Value = " Для начала, нам нужно достичь\nРадужного Камнерабля!"
$ws.Range("E17").Value = " Äìÿ îàœàìà, îàí îôçîï äïòóéœû\nÑàäôçîïãï Ëàíîåñàáìÿ!"
$ws.Range("A18").Value = "SCRIPT/D01P11A/um2401.ssb"
$ws.Range("B18").Value = "308, 112"
$ws.Range("C18").Value = " We\'ll leave when you\'re ready!"
$ws.Range("D18").Value = " Мы покинем пляж сразу же, как\nвы подготовитесь!"
$ws.Range("E18").Value = " Íú ðïëéîåí ðìÿç òñàèô çå, ëàë\nâú ðïäãïóïâéóåòû!"
$ws.Range("B19").Value = "311, 115"
$ws.Range("C19").Value = " We\'re going off to the\n[CS:P]Hidden Land[CR]!"
$ws.Range("D19").Value = " Мы отправимся в [CS:P]Сокрытые Земли[CR]!"
$ws.Range("E19").Value = " Íú ïóðñàâéíòÿ â [CS:P]Òïëñúóúå Èåíìé[CR]!"
$ws.Range("A20").Value = "SCRIPT/D27P11A/um2404.ssb"
$ws.Range("B20").Value = "159, 90"
$ws.Range("C20").Value = " So the priority is finding the\nRainbow Stoneship.[K] If we find it, we can get\nto [CS:P]Temporal Tower[CR]."
$ws.Range("D20").Value = " Наша основная задача - найти\nРадужный Камнерабль.[K] Если мы отыщем его,\nто сможем попасть в [CS:P]Темпоральную Башню[CR]."
$ws.Range("E20").Value = " Îàšà ïòîïâîàÿ èàäàœà - îàêóé\nÑàäôçîúê Ëàíîåñàáìû.[K] Åòìé íú ïóúþåí åãï,\nóï òíïçåí ðïðàòóû â [CS:P]Óåíðïñàìûîôý Áàšîý[CR]."
$ws.Range("B21").Value = "162, 93"
$ws.Range("C21").Value = " We\'ll leave when you\'re ready!"
$ws.Range("D21").Value = " Мы уйдём сразу же, как вы\nподготовитесь!"
$ws.Range("E21").Value = " Íú ôêäæí òñàèô çå, ëàë âú\nðïäãïóïâéóåòû!"
$ws.Range("A22").Value = "SCRIPT/D28P21A/um2401.ssb"
$ws.Range("B22").Value = 203
$ws.Range("C22").Value = " According to [CS:N]Lapras[CR], the\n[CS:P]Old Ruins[CR] lie beyond this dungeon…"
$ws.Range("D22").Value = " Судя по тому, что говорил\n[CS:N]Лапрас[CR], [CS:P]Старые Руины[CR] находятся прямо за\nэтим подземельем..."
$ws.Range("E22").Value = " Òôäÿ ðï óïíô, œóï ãïâïñéì\n[CS:N]Ìàðñàò[CR], [CS:P]Òóàñúå Ñôéîú[CR] îàöïäÿóòÿ ðñÿíï èà\nüóéí ðïäèåíåìûåí..."
$ws.Range("B23").Value = 206
$ws.Range("C23").Value = " That\'s where we should find the\nRainbow Stoneship."
$ws.Range("D23").Value = " Там мы сможем найти Радужный\nКамнерабль."
$ws.Range("E23").Value = " Óàí íú òíïçåí îàêóé Ñàäôçîúê\nËàíîåñàáìû."
$ws.Range("B24").Value = 209
$ws.Range("C24").Value = " If we find the Rainbow\nStoneship, we\'ll take that to [CS:P]Temporal Tower[CR]."
$ws.Range("D24").Value = " Если мы найдём Радужный\nКамнерабль, то сможем попасть в\n[CS:P]Темпоральную Башню[CR]."
$ws.Range("E24").Value = " Åòìé íú îàêäæí Ñàäôçîúê\nËàíîåñàáìû, óï òíïçåí ðïðàòóû â\n[CS:P]Óåíðïñàìûîôý Áàšîý[CR]."
$ws.Range("A25").Value = "SCRIPT/D28P21A/um2402.ssb"
$ws.Range("B25").Value = 134
$ws.Range("C25").Value = " We\'ve made it this far."
$ws.Range("D25").Value = " Мы уже далеко зашли."
$ws.Range("E25").Value = " Íú ôçå äàìåëï èàšìé."
$ws.Range("B26").Value = 167
$ws.Range("C26").Value = " We\'ll get through this,\nno matter what."
$ws.Range("D26").Value = " Так или иначе, но мы продвинемся\nдальше."
$ws.Range("E26").Value = " Óàë éìé éîàœå, îï íú ðñïäâéîåíòÿ\näàìûšå. "
$ws.Range("A27").Value = "SCRIPT/D28P21A/um2403.ssb"
$ws.Range("B27").Value = 65
$ws.Range("C27").Value = " It doesn\'t matter if [CS:N]Dusknoir[CR] is\nthere or not."
$ws.Range("D27").Value = " Не имеет значения, засел\n[CS:N]Даскнуар[CR] в засаде или нет."
$ws.Range("E27").Value = " Îå éíååó èîàœåîéÿ, èàòåì\n[CS:N]Äàòëîôàñ[CR] â èàòàäå éìé îåó."
$ws.Range("B28").Value = 68
$ws.Range("C28").Value = " We\'ll get through next time."
$ws.Range("D28").Value = " На этот раз мы прорвёмся."
$ws.Range("E28").Value = " Îà üóïó ñàè íú ðñïñâæíòÿ."
$ws.Range("A29").Value = "SCRIPT/D27P11A/um2404.ssb "
$ws.Range("B29").Value = 20
$ws.Range("C29").Value = " I don\'t care if [CS:N]Dusknoir[CR] is\nwaiting to ambush us."
$ws.Range("D29").Value = " Мне всё равно, засел\n[CS:N]Даскнуар[CR] в засаде или нет."
$ws.Range("E29").Value = " Íîå âòæ ñàâîï, èàòåì\n[CS:N]Äàòëîôàñ[CR] â èàòàäå éìé îåó."
$ws.Range("B30").Value = 23
$ws.Range("C30").Value = " We have to get through and get\non the Rainbow Stoneship."
$ws.Range("D30").Value = " На этот раз мы прорвёмся к\nРадужному Камнераблю."
$ws.Range("E30").Value = " Îà üóïó ñàè íú ðñïñâæíòÿ ë\nÑàäôçîïíô Ëàíîåñàáìý."
$ws.Range("A31").Value = "SCRIPT/D01P11A/um2401.ssb"
$ws.Range("B31").Value = 43
$ws.Range("C31").Value = " We\'ll leave when you\'re ready!"
$ws.Range("D31").Value = " Как только вы подготовитесь,\nмы тут же покинем пляж!"
$ws.Range("E31").Value = " Ëàë óïìûëï âú ðïäãïóïâéóåòû,\níú óôó çå ðïëéîåí ðìÿç!"
$ws.Range("B32").Value = 46
$ws.Range("C32").Value = " We\'re going off to the\n[CS:P]Hidden Land[CR]!"
$ws.Range("D32").Value = " Мы отправимся в [CS:P]Сокрытые Земли[CR]!"
$ws.Range("E32").Value = " Íú ïóðñàâéíòÿ â [CS:P]Òïëñúóúå Èåíìé[CR]!"

# --- Clear clipboard marching ants / selection / view state ---
$excel.CutCopyMode = $false
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("C28").Select()
